# Final push, whole project and outputs
#
# Updates the BOM report's "generated on" metadata (report date/time) and
# bumps two component quantities (H12, H13) on the "Part List Report"
# sheet, mirroring the same values on the "Project Information" sheet.
# A leading apostrophe is used for the textual date/time/quantity values
# so Excel keeps storing them as plain text (matching the original
# quotePrefix-styled cells) instead of re-typing them as numbers/dates.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Part List Report")
$ws2 = $wb.Worksheets.Item("Project Information")

# --- "Part List Report" sheet -------------------------------------------------
# Report Date / Report Time shown next to the TODAY()/NOW() cells.
$ws1.Range("D7").Value = "'23/10/2023"
$ws1.Range("E7").Value = "'15:17"

# Quantities for the two SMD test point line items went from 1 -> 2 each.
$ws1.Range("H12").Value = 2
$ws1.Range("H13").Value = 2

# --- "Project Information" sheet ----------------------------------------------
# Total Quantity (sum of H10:H13 on the other sheet) 8 -> 10.
$ws2.Range("B7").Value = "'10"
# Report Time / Report Date / Report Date & Time.
$ws2.Range("B8").Value = "'15:17"
$ws2.Range("B9").Value = "'23/10/2023"
$ws2.Range("B10").Value = "'23/10/2023 15:17"
